$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update coefficients / RHS values (fixing simplex RHS bugs)
$ws.Range("A2").Value = 4
$ws.Range("B2").Value = 3

$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 40

$ws.Range("A4").Value = 6
$ws.Range("B4").Value = -2
$ws.Range("C4").Value = 30

# Clear out the now-unused 5th constraint row, keeping styles intact
$ws.Range("A5").ClearContents()
$ws.Range("B5").ClearContents()
$ws.Range("C5").ClearContents()
$ws.Range("E5").ClearContents()
